$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the data refresh diff (price/volume values, and two coin-row swaps)
$updates = @(
    @('D2', '63.582.36'),
    @('E2', '  +4.65%  '),
    @('D3', '3.424.91'),
    @('E3', '  +5.89%  '),
    @('E4', '  -0.02%  '),
    @('D5', '576.12'),
    @('E5', '  +6.39%  '),
    @('E6', '  +6.40%  '),
    @('D7', '0.999'),
    @('E7', '  -0.02%  '),
    @('D8', '3.427.33'),
    @('E8', '  +5.58%  '),
    @('E9', '  +1.35%  '),
    @('E10', '  +2.71%  '),
    @('E11', '  +7.48%  '),
    @('D12', '0.438'),
    @('E12', '  +0.25%  '),
    @('D13', '4.013.63'),
    @('E13', '  +5.97%  '),
    @('E14', '  -0.69%  '),
    @('E15', '  +7.03%  '),
    @('D16', '27.35'),
    @('E16', '  +4.61%  '),
    @('D17', '63.682.81'),
    @('E17', '  +4.81%  '),
    @('D18', '3.423.69'),
    @('E18', '  +5.93%  '),
    @('D19', '6.44'),
    @('E19', '  +1.82%  '),
    @('D20', '14.31'),
    @('E20', '  +7.19%  '),
    @('E21', '  +1.79%  '),
    @('D22', '391.89'),
    @('E22', '  +3.84%  '),
    @('E23', '  -0.18%  '),
    @('E24', '  +2.04%  '),
    @('D25', '71.97'),
    @('E25', '  +2.74%  '),
    @('E26', '  +19.10%  '),
    @('D27', '9.54'),
    @('E27', '  +9.94%  '),
    @('E28', '  +5.09%  '),
    @('E29', '  -0.11%  '),
    @('E30', '  +7.72%  '),
    @('E31', '  +12.70%  '),
    @('E32', '  +6.47%  '),
    @('D33', '5.82'),
    @('E33', '  +8.27%  '),
    @('D34', '23.57'),
    @('E34', '  +4.36%  '),
    @('E35', '  -0.12%  '),
    @('E36', '  +3.40%  '),
    @('D37', '1.50'),
    @('E37', '  +4.40%  '),
    @('D38', '158.18'),
    @('E38', '  +0.03%  '),
    @('D39', '28.16'),
    @('E39', '  +6.41%  '),
    @('D40', '0.0785'),
    @('E40', '  +9.81%  '),
    @('D41', '1.86'),
    @('E41', '  +8.04%  '),
    @('D42', '2.865.74'),
    @('E42', '  +2.09%  '),
    @('D43', '0.0319'),
    @('E43', '  +1.78%  '),
    @('B44', 'OKB'),
    @('C44', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
    @('D44', '41.93'),
    @('E44', '  +4.87%  '),
    @('B45', 'Mantle'),
    @('C45', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D45', '0.768'),
    @('E45', '  +5.98%  '),
    @('D46', '4.39'),
    @('E46', '  +2.59%  '),
    @('E47', '  +9.27%  '),
    @('D48', '3.468.87'),
    @('E48', '  +5.98%  '),
    @('E49', '  +7.18%  '),
    @('B50', 'dogwifhat'),
    @('C50', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'),
    @('D50', '2.11'),
    @('E50', '  +23.05%  '),
    @('B51', 'Cosmos'),
    @('C51', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @('D51', '6.38'),
    @('E51', '  +2.80%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "0.999") are not
    # reinterpreted as numbers by Excel, matching the original inlineStr text cells.
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
    $rng.Style = "Normal"
}
